# Updated symbol list on Thu Jan 26 20:49:34 UTC 2023 with GitHub Actions
#
# This script rewrites the "Price" (D) and "Volume(1h)" (E) columns with
# refreshed quotes, and fixes a swapped pair of rows (20/21: MCDex <->
# ProBitToken) so that the coin/link/price/volume line back up correctly.
#
# Price/Volume cells hold numeric-looking text (e.g. "306.65", "1.11%") but
# must stay plain text, matching the original workbook. Excel's COM layer
# auto-coerces such strings into real numbers on a plain `.Value =`
# assignment, so for those cells we momentarily force a text number format,
# assign the value, then restore the "Normal" style so no stray formatting
# is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: named parameter binding (-CellRef "...") is unreliable on this
# engine, so this helper is always invoked positionally.
function Set-TextValue {
    param([string]$CellRef, [string]$Text)
    $cell = $ws.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# Row -> new Price (D) / Volume(1h) (E) values
$deChanges = @{
    2  = @("306.65", "1.11%")
    3  = @("35.97", "0.27%")
    4  = @("5.000", "-1.28%")
    5  = @("0.08087", "0.28%")
    6  = @("1.925", "-0.51%")
    7  = @("4.144", "2.21%")
    8  = @("7.879", "1.17%")
    9  = @("0.9297", "-0.03%")
    10 = @("0.1255", "-17.20%")
    11 = @("0.1901", "-0.29%")
    12 = @("0.09233", "2.80%")
    13 = @("0.03510", "1.48%")
    14 = @("0.09904", "0.62%")
    15 = @("0.001417", "1.46%")
    16 = @("0.006604", "14.95%")
    17 = @("3.613", "2.10%")
    19 = @("0.3438", "-0.19%")
    22 = @("0.2533", "5.80%")
    23 = @("0.04408", "-1.99%")
    24 = @("0.001233", "2.03%")
    25 = @("0.004730", "-1.63%")
    26 = @("0.0001299", "5.87%")
    27 = @("0.0003126", "3.69%")
    39 = @("0.01966", "4.60%")
    40 = @("0.05280", "10.11%")
    41 = @("0.007547", "3.61%")
    42 = @("0.01014", "-4.34%")
    43 = @("0.1371", "2.00%")
    44 = @("0.002098", "-0.31%")
    45 = @("0.01066", "9.69%")
    46 = @("0.00006383", "2.65%")
    48 = @("63.57", "-1.70%")
    49 = @("0.001658", "-0.11%")
    50 = @("0.00002099", "0.20%")
    51 = @("0.0001999", "0.20%")
}

foreach ($row in $deChanges.Keys) {
    $vals = $deChanges[$row]
    Set-TextValue ("D{0}" -f $row) $vals[0]
    Set-TextValue ("E{0}" -f $row) $vals[1]
}

# Row 47 only has a Volume(1h) change (its Price is a non-numeric "--").
Set-TextValue "E47" "0.20%"

# Rows 20/21 were re-ordered upstream: MCDex and ProBitToken swapped places,
# each keeping the rank (column A), date (F) and hour (G) of its row but
# bringing its own coin/link/price/volume along.
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D20" "0.1329"
Set-TextValue "E20" "2.33%"

$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D21" "5.192"
Set-TextValue "E21" "3.06%"
